# Refactor tournament registration data:
# - Every account row (2-32) now shows it is registered for the
#   "10.000 TL Texas Ücretsiz Turnuva" tournament (previously only row 2 had
#   a (different) tournament name, the rest were blank).
# - Rows 3, 13 and 31 are now flagged as registered (registered = 1).
# - Rows 2, 4 and 13 no longer carry stale balance/last-check-time data for
#   the old tournament, so those columns are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tournament = "10.000 TL Texas Ücretsiz Turnuva"

# Update the registered_tournament column (E) for every data row (2-32).
for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 5).Value = $tournament
}

# Mark these rows as registered (column D = 1).
$registeredRows = 3, 13, 31
foreach ($row in $registeredRows) {
    $ws.Cells.Item($row, 4).Value = 1
}

# Clear out the now-stale balance / last-check-time columns (F:I) for the
# rows whose tournament info just changed.
$clearRows = 2, 4, 13
foreach ($row in $clearRows) {
    $ws.Range($ws.Cells.Item($row, 6), $ws.Cells.Item($row, 9)).ClearContents()
}
